$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76 (anchor G=12602)
$ws.Range("H76").Value = 2957.5
$ws.Range("I76").Value = 2957.5
$ws.Range("K76").Value = 2957.5
$ws.Range("M76").Value = -2642.5
# Row 79 (anchor G=12602)
$ws.Range("H79").Value = 2957.5
$ws.Range("I79").Value = 2957.5
$ws.Range("K79").Value = 2957.5
$ws.Range("M79").Value = -1865.5
# Row 111 (anchor G=27768)
$ws.Range("H111").Value = 2664.9583
$ws.Range("I111").Value = 2311.2
$ws.Range("J111").Value = 3254.5557
$ws.Range("K111").Value = 6933.599999999999
$ws.Range("L111").Value = 9763.667099999999
$ws.Range("M111").Value = -3866.599999999999
$ws.Range("N111").Value = -15897.6671
# Row 132 (anchor G=44049)
$ws.Range("H132").Value = 21789732
$ws.Range("I132").Value = 2316131.5
$ws.Range("J132").Value = 333367330
$ws.Range("K132").Value = 6948394.5
$ws.Range("L132").Value = 1000101990
$ws.Range("M132").Value = -6945864.5
$ws.Range("N132").Value = -1000107050
# Row 137 (anchor G=44013)
$ws.Range("H137").Value = 1430.3959
$ws.Range("I137").Value = 1072.5
$ws.Range("J137").Value = 2146.1875
$ws.Range("K137").Value = 3217.5
$ws.Range("L137").Value = 6438.5625
$ws.Range("M137").Value = -667.5
$ws.Range("N137").Value = -11538.5625
# Row 141 (anchor G=44161)
$ws.Range("H141").Value = 997.8
$ws.Range("I141").Value = 609.2143
$ws.Range("K141").Value = 1827.6429
$ws.Range("M141").Value = 3352.3571

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (anchor G=44147)
$ws.Range("H32").Value = 4552.3125
$ws.Range("I32").Value = 2995.9268
$ws.Range("J32").Value = 13668.286
$ws.Range("K32").Value = 2995.9268
$ws.Range("L32").Value = 13668.286
$ws.Range("M32").Value = -2708.9268
$ws.Range("N32").Value = -14242.286
# Row 36 (anchor G=3068)
$ws.Range("H36").Value = 5993.3335
$ws.Range("I36").Value = 5993.3335
$ws.Range("K36").Value = 5993.3335
$ws.Range("M36").Value = -5647.3335
# Row 132 (anchor G=43997)
$ws.Range("H132").Value = 1781.0652
$ws.Range("I132").Value = 1775.5
$ws.Range("J132").Value = 1793.7858
$ws.Range("K132").Value = 5326.5
$ws.Range("L132").Value = 5381.357400000001
$ws.Range("M132").Value = -2796.5
$ws.Range("N132").Value = -10441.3574
# Row 140 (anchor G=42496)
$ws.Range("H140").Value = 28849
$ws.Range("J140").Value = 28849
$ws.Range("L140").Value = 28849
$ws.Range("N140").Value = -39209

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 33 (anchor G=1625)
$ws.Range("H33").Value = 50000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 50000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 50000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -50672
# Row 55 (anchor G=27151)
$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5
# Row 99 (anchor G=19943)
$ws.Range("H99").Value = 1747.52
$ws.Range("I99").Value = 1666.0555
$ws.Range("J99").Value = 1957
$ws.Range("K99").Value = 1666.0555
$ws.Range("L99").Value = 1957
$ws.Range("M99").Value = -168.0554999999999
$ws.Range("N99").Value = -4953

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25 (anchor G=1895)
$ws.Range("H25").Value = 4876.15
$ws.Range("I25").Value = 2032.6875
$ws.Range("K25").Value = 2032.6875
$ws.Range("M25").Value = -1858.6875
# Row 39 (anchor G=1915)
$ws.Range("H39").Value = 2040.8
$ws.Range("I39").Value = 2040.8
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2040.8
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -1649.8
$ws.Range("N39").ClearContents()
# Row 49 (anchor G=1915)
$ws.Range("H49").Value = 2040.8
$ws.Range("I49").Value = 2040.8
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2040.8
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -1858.8
$ws.Range("N49").ClearContents()
# Row 104 (anchor G=19749)
$ws.Range("H104").Value = 28000
$ws.Range("J104").Value = 28000
$ws.Range("L104").Value = 28000
$ws.Range("N104").Value = -33242
# Row 134 (anchor G=44020)
$ws.Range("H134").Value = 25001910
$ws.Range("I134").Value = 2704532.8
$ws.Range("J134").Value = 142859470
$ws.Range("K134").Value = 8113598.399999999
$ws.Range("L134").Value = 428578410
$ws.Range("M134").Value = -8111063.399999999
$ws.Range("N134").Value = -428583480
# Row 141 (anchor G=43345)
$ws.Range("H141").Value = 58259.8
$ws.Range("J141").Value = 61366.89
$ws.Range("L141").Value = 61366.89
$ws.Range("N141").Value = -71726.89

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 109 (anchor G=27854)
$ws.Range("H109").Value = 1757.8636
$ws.Range("I109").Value = 912.5714
$ws.Range("J109").Value = 3237.125
$ws.Range("K109").Value = 2737.7142
$ws.Range("L109").Value = 9711.375
$ws.Range("M109").Value = -1697.7142
$ws.Range("N109").Value = -11791.375
# Row 131 (anchor G=36060)
$ws.Range("H131").Value = 909.86
$ws.Range("I131").Value = 818.5714
$ws.Range("J131").Value = 916.7311999999999
$ws.Range("K131").Value = 2455.7142
$ws.Range("L131").Value = 2750.1936
$ws.Range("M131").Value = 2584.2858
$ws.Range("N131").Value = -12830.1936
# Row 139 (anchor G=44102)
$ws.Range("H139").Value = 2202.75
$ws.Range("I139").Value = 1731.0667
$ws.Range("J139").Value = 2988.889
$ws.Range("K139").Value = 5193.2001
$ws.Range("L139").Value = 8966.667000000001
$ws.Range("M139").Value = -53.20010000000002
$ws.Range("N139").Value = -19246.667

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107 (anchor G=27802)
$ws.Range("H107").Value = 1142.05
$ws.Range("I107").Value = 1026.0667
$ws.Range("J107").Value = 1490
$ws.Range("K107").Value = 1026.0667
$ws.Range("L107").Value = 1490
$ws.Range("M107").Value = 893.9332999999999
$ws.Range("N107").Value = -5330

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (anchor G=36249)
$ws.Range("H7").Value = 2833
$ws.Range("I7").Value = 2750
$ws.Range("K7").Value = 2750
$ws.Range("M7").Value = -2638
# Row 40 (anchor G=36248)
$ws.Range("H40").Value = 2940.4443
$ws.Range("I40").Value = 2934.875
$ws.Range("J40").Value = 2985
$ws.Range("K40").Value = 2934.875
$ws.Range("L40").Value = 2985
$ws.Range("M40").Value = -2798.875
$ws.Range("N40").Value = -3257
# Row 68 (anchor G=12563)
$ws.Range("H68").Value = 6520.36
$ws.Range("I68").Value = 14336.25
$ws.Range("J68").Value = 2842.2942
$ws.Range("K68").Value = 14336.25
$ws.Range("L68").Value = 2842.2942
$ws.Range("M68").Value = -13587.25
$ws.Range("N68").Value = -4340.2942
# Row 71 (anchor G=12563)
$ws.Range("H71").Value = 6520.36
$ws.Range("I71").Value = 14336.25
$ws.Range("J71").Value = 2842.2942
$ws.Range("K71").Value = 71681.25
$ws.Range("L71").Value = 14211.471
$ws.Range("M71").Value = -67937.25
$ws.Range("N71").Value = -21699.471
# Row 122 (anchor G=36247)
$ws.Range("H122").Value = 10600
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -17900.0005
# Row 126 (anchor G=36249)
$ws.Range("H126").Value = 2833
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780
# Row 132 (anchor G=44058)
$ws.Range("H132").Value = 3682.1428
$ws.Range("I132").Value = 4770.4
$ws.Range("J132").Value = 2231.1333
$ws.Range("K132").Value = 14311.2
$ws.Range("L132").Value = 6693.3999
$ws.Range("M132").Value = -11781.2
$ws.Range("N132").Value = -11753.3999
# Row 136 (anchor G=44060)
$ws.Range("H136").Value = 1685.3334
$ws.Range("I136").Value = 1123.6
$ws.Range("K136").Value = 3370.8
$ws.Range("M136").Value = -820.7999999999997

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5 (anchor G=3515)
$ws.Range("H5").Value = 30000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 30000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 30000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -30224

